$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.120.92"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.95%  "
$ws.Range("D3").Value = "'2.000.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.70%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'242.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.70%  "
$ws.Range("D6").Value = "'0.605"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.92%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'55.22"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.75%  "
$ws.Range("D9").Value = "'0.374"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.35%  "
$ws.Range("D10").Value = "'58.40"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.92%  "
$ws.Range("D11").Value = "'0.0757"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.95%  "
$ws.Range("D12").Value = "'0.0980"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.27%  "
$ws.Range("D13").Value = "'2.290.14"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.00%  "
$ws.Range("D14").Value = "'14.15"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.34%  "
$ws.Range("D15").Value = "'20.89"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.63%  "
$ws.Range("D16").Value = "'0.761"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -7.54%  "
$ws.Range("D17").Value = "'5.06"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.90%  "
$ws.Range("D18").Value = "'1.999.18"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.05%  "
$ws.Range("D19").Value = "'37.008.49"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.38%  "
$ws.Range("D20").Value = "'68.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.58%  "
$ws.Range("D21").Value = "'0.0₃0812"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.19%  "
$ws.Range("D22").Value = "'229.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("D23").Value = "'5.01"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.96%  "
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").Value = "'2.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -8.64%  "
$ws.Range("E26").Value = "  -0.21%  "
$ws.Range("D27").Value = "'162.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.03%  "
$ws.Range("D28").Value = "'8.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.16%  "
$ws.Range("D29").Value = "'19.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.82%  "
$ws.Range("D30").Value = "'0.124"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -9.15%  "
$ws.Range("D31").Value = "'1.32"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.62%  "
$ws.Range("E32").Value = "  -2.75%  "
$ws.Range("D33").Value = "'4.45"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.08%  "
$ws.Range("D34").Value = "'0.0613"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -8.21%  "
$ws.Range("D35").Value = "'4.24"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.58%  "
$ws.Range("D36").Value = "'2.35"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.86%  "
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("E38").Value = "  -1.63%  "
$ws.Range("D39").Value = "'3.33"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.87%  "
$ws.Range("D40").Value = "'5.23"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.81%  "
$ws.Range("E41").Value = "  +2.51%  "
$ws.Range("D42").Value = "'1.439.69"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.47%  "
$ws.Range("D43").Value = "'0.0204"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.71%  "
$ws.Range("E44").Value = "  -5.81%  "
$ws.Range("D45").Value = "'0.0883"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -8.80%  "
$ws.Range("D46").Value = "'88.67"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.61%  "
$ws.Range("D47").Value = "'15.28"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.31%  "
$ws.Range("D48").Value = "'1.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.07%  "
$ws.Range("E49").Value = "  +0.72%  "
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").Value = "'6.68"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -9.51%  "
$ws.Range("B51").Value = "FTXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D51").Value = "'3.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +12.46%  "
